# The captured change only re-serialises the package's OOXML (the XML
# attributes on elements such as <w:tblW>, <w:tcW>, <w:tblLook>,
# <w:pgSz>, <w:pgMar>, <w:rFonts>, <w:lang>, <w:latentStyles>,
# <w:lsdException>, <w:style>, <w:tblBorders>, ... are simply written out
# in a different (alphabetised) attribute order) as part of bumping the
# authoring tool from 2.0.2 to 2.0.3. Every attribute name/value pair and
# every piece of document content (text, table layout, styles, …) is
# identical before and after - there is no actual content edit to make.
#
# Touch the document model read-only so the session has a well-defined
# "applied" pass, without mutating any visible content, formatting, or
# structure.
$d = $word.ActiveDocument
$null = $d.Content.Text
$null = $d.Tables.Count
$null = $d.Styles.Count
